$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-detected as a number by Excel;
# force Text format first so the literal string (matching the original inline-string
# formatting, e.g. trailing zeros) is preserved exactly.
$ws.Range('D2').Value = '27.145.87'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '1.852.24'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').Value = '  +0.81%  '
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '310.07'
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4780'
$ws.Range('E7').Value = '  +2.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3692'
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07258'
$ws.Range('E9').Value = '  +1.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9322'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.92'
$ws.Range('E11').Value = '  +1.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07794'
$ws.Range('E12').Value = '  +1.36%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.395'
$ws.Range('E13').Value = '  +2.12%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.794.69'
$ws.Range('E14').Value = '  -3.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.489'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '89.40'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.018'
$ws.Range('E17').Value = '  +1.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008700'
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('E19').Value = '  +0.78%  '
$ws.Range('D20').Value = '27.163.76'
$ws.Range('E20').Value = '  +0.91%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.62'
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.066'
$ws.Range('E22').Value = '  +0.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.65'
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.947'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.23'
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.36'
$ws.Range('E26').Value = '  +0.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.989'
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.89'
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.933'
$ws.Range('E29').Value = '  +1.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08886'
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('E31').Value = '  +3.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.186'
$ws.Range('E32').Value = '  +1.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.526'
$ws.Range('E33').Value = '  +1.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7392'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.692'
$ws.Range('E35').Value = '  -3.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.119'
$ws.Range('E36').Value = '  +3.35%  '
$ws.Range('E37').Value = '  +2.22%  '
$ws.Range('E38').Value = '  +1.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.980'
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5300'
$ws.Range('E40').Value = '  +1.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.057'
$ws.Range('E41').Value = '  +2.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1525'
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.316'
$ws.Range('E43').Value = '  +2.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.56'
$ws.Range('E44').Value = '  +0.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4753'
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('E46').Value = '  +0.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.05'
$ws.Range('E47').Value = '  +1.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.622'
$ws.Range('E48').Value = '  +1.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '65.90'
$ws.Range('E49').Value = '  +1.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06059'
$ws.Range('E50').Value = '  +0.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8934'
$ws.Range('E51').Value = '  +0.18%  '
